$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 - headers
$ws.Range("A1").Value = "code"
$ws.Range("B1").Value = "name"
$ws.Range("C1").Value = "descr"
$ws.Range("D1").Value = "lang_code"
$ws.Range("E1").Value = "is_active"
$ws.Range("F1").Value = "cr_by"
$ws.Range("G1").Value = "cr_dtimes"
$ws.Range("H1").Value = "upd_by"
$ws.Range("I1").Value = "upd_dtimes"
$ws.Range("J1").Value = "is_deleted"
$ws.Range("K1").Value = "del_dtimes"

# Row 2 - data
$ws.Range("A2").Value = "ADT"
$ws.Range("B2").Value = "Activation status"
$ws.Range("C2").Value = "Activation status"
$ws.Range("D2").Value = "fra"
$ws.Range("E2").Value = $true
$ws.Range("F2").Value = "superadmin"
$ws.Range("G2").Value = 45079.577392824074
$ws.Range("G2").NumberFormat = "mm:ss.0"
$ws.Range("H2").Value = "NULL"
$ws.Range("I2").Value = "NULL"
$ws.Range("J2").Value = $false
$ws.Range("K2").Value = "NULL"

# Selection matches target
[void]$ws.Range("D8").Select()
